$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the two "pro003" / "pro027" test-case rows to use btree_-prefixed
# user names (the commit renames them so they don't collide with other
# engines' test fixtures). ---

# Row 4 : create/show user 'pro003'
$ws.Range("F4").Value = "btree_pro003"
$ws.Range("H4").Value = "create user 'btree_pro003'@'172.20.3.15' identified by 'abc123'"
$ws.Range("I4").Value = "show create user 'btree_pro003'@'172.20.3.15'"

# Row 28 : create/show user 'pro027'
$ws.Range("F28").Value = "btree_pro027"
$ws.Range("H28").Value = "create user 'btree_pro027' identified by 'abc123' require ssl"
$ws.Range("I28").Value = "show create user 'btree_pro027'"

# --- View/formatting tweaks captured in the diff ---

# Column F a bit wider to fit "btree_pro003" / "btree_pro027"
$ws.Columns("F").ColumnWidth = 13.875

# Scroll back to the top of the sheet and select F3 (previously topLeftCell
# was A7 with H29 selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F3").Select()
